# The workbook tracks weekly Mango price records for "Femacal de La Calera".
# This edit adds one new weekly record. In the canonical row order the new
# record lands at row 286, which pushes every following record down by one
# row (so the former last row, 407, becomes row 408).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 286; Excel shifts rows 286-407 down to 287-408,
# carrying their values/formatting with them (matches the data-only shift
# seen for every row after the insertion point).
$ws.Rows("286:286").Insert()

# Populate the newly-inserted row with the new weekly record. The
# "descriptive" columns (mercado/region/producto/etc.) repeat the same
# constants used throughout this sheet.
$ws.Range("A286").Value = 3
$ws.Range("B286").Value = "Femacal de La Calera"
$ws.Range("C286").Value = "Coquimbo"
$ws.Range("D286").Value = 44704
$ws.Range("E286").Value = 5
$ws.Range("F286").Value = "Fruta"
$ws.Range("G286").Value = 100108
$ws.Range("H286").Value = "Tropicales y subtropicales"
$ws.Range("I286").Value = 100108002
$ws.Range("J286").Value = "Mango"
$ws.Range("K286").Value = "Sin especificar"
$ws.Range("L286").Value = "Primera"
$ws.Range("M286").Value = 228
$ws.Range("N286").Value = 9000
$ws.Range("O286").Value = 9000
$ws.Range("P286").Value = 9000
$ws.Range("Q286").Value = "`$/bandeja 4 kilos"
$ws.Range("R286").Value = "Brasil"
$ws.Range("S286").Value = 2250
$ws.Range("T286").Value = 4
